$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect to edit, then restore protection after.
$ws.Unprotect()

# Update the confidential disclaimer text (date changed from 2021-04-06 to 2021-04-08)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1027368546321515
$ws.Range("E2").Value = 0.01071347162069758

$ws.Range("D3").Value = 0.1071865300091939
$ws.Range("E3").Value = 0.0157215721572157

$ws.Range("D4").Value = 0.1170816607538862
$ws.Range("E4").Value = 0.004014452027298221

$ws.Range("D5").Value = 0.1369176047635423
$ws.Range("E5").Value = 0.003358178392988576

$ws.Range("D6").Value = 0.1334424348460063
$ws.Range("E6").Value = 0.001847640704945919

$ws.Range("D7").Value = 0.1447934691609345
$ws.Range("E7").Value = -0.0009567546880979716

$ws.Range("D8").Value = 0.1295734387248316
$ws.Range("E8").Value = 0.008489461358313966

$ws.Range("D9").Value = 0.1282680071094535
$ws.Range("E9").Value = 0.01383596369286111

$ws.Range("D10").Value = 0.9999999999999998
$ws.Range("E10").Value = 0.006698363629035198

# Restore sheet protection (matches original protection settings)
$ws.Protect($null, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true, $true, $false, $false, $true)
